$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the (unstyled) format of a plain data cell so we can force text
# storage for numeric-looking strings without altering cell styling.
$plainStyle = $ws.Range("B2").Style

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.100.61"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.796.10"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +2.02%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.49"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +1.22%  "

# Row 6
$ws.Range("E6").Value = "  +0.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +23.34%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3724"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +10.68%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07680"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +6.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.146"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +1.37%  "

# Row 12
$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.53"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  -0.90%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.398"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +3.52%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.479"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +3.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.796.98"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +2.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  +3.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06751"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +2.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.36"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("E20").Value = "  +0.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.51"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +2.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.446"
$ws.Range("D22").Style = $plainStyle

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.124.74"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.97"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +2.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.407"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +0.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.86"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +4.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.406"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +3.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.30"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -1.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.004.34"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  +2.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.38"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +1.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.262"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -0.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.055"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +0.97%  "

# Row 33
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09642"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +9.94%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.945"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +2.11%  "

# Row 35
$ws.Range("E35").Value = "  +2.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2229"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +5.26%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.20"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06390"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +2.77%  "

# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6722"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +0.93%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.274"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +2.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.237"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +1.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.488"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +2.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.082"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +0.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.27"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +3.96%  "

# Row 45
$ws.Range("E45").Value = "  +0.40%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6159"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +1.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.864"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +1.40%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.99"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +0.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.068"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +2.68%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.181"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  -1.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07129"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -1.14%  "
